$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 654, pushing the existing rows 654-700 down to 655-701
$ws.Rows("654:654").Insert()

# Populate the newly inserted row 654 with the new weekly price record
$ws.Range("A654").Value = 5
$ws.Range("B654").Value = "Macroferia Regional de Talca"
$ws.Range("C654").Value = "Maule"
$ws.Range("D654").Value = 45265
$ws.Range("E654").Value = 7
$ws.Range("F654").Value = 100114014
$ws.Range("G654").Value = "Betarraga"
$ws.Range("H654").Value = "Sin especificar"
$ws.Range("I654").Value = "Primera"
$ws.Range("J654").Value = 5000
$ws.Range("K654").Value = 700
$ws.Range("L654").Value = 700
$ws.Range("M654").Value = 700
$ws.Range("N654").Value = "$/paquete 5 unidades"
$ws.Range("O654").Value = "Región del Maule"
$ws.Range("P654").Value = 140
$ws.Range("Q654").Value = 5
$ws.Range("R654").Value = "Hortaliza"
